$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart (item id 7015)
$ws.Cells.Item(19, 8).Value = 1180
$ws.Cells.Item(19, 9).Value = 1125
$ws.Cells.Item(19, 10).Value = 1207.5
$ws.Cells.Item(19, 11).Value = 1125
$ws.Cells.Item(19, 12).Value = 1207.5
$ws.Cells.Item(19, 13).Value = -950
$ws.Cells.Item(19, 14).Value = -1557.5

# Row 32: Automata for the People (item id 5484)
$ws.Cells.Item(32, 8).Value = 8783.15
$ws.Cells.Item(32, 9).Value = 12875
$ws.Cells.Item(32, 10).Value = 7029.5
$ws.Cells.Item(32, 11).Value = 12875
$ws.Cells.Item(32, 12).Value = 7029.5
$ws.Cells.Item(32, 13).Value = -12549
$ws.Cells.Item(32, 14).Value = -7681.5

# Row 40: Stuck in the Moment (item id 5505)
$ws.Cells.Item(40, 8).Value = 1545.2941
$ws.Cells.Item(40, 9).Value = 1569.2307
$ws.Cells.Item(40, 10).Value = 1467.5
$ws.Cells.Item(40, 11).Value = 1569.2307
$ws.Cells.Item(40, 12).Value = 1467.5
$ws.Cells.Item(40, 13).Value = -1394.2307
$ws.Cells.Item(40, 14).Value = -1817.5

# Row 99: Rumor Has It (item id 19883)
$ws.Cells.Item(99, 8).Value = 6197.25
$ws.Cells.Item(99, 9).Value = 2808.5
$ws.Cells.Item(99, 10).Value = 7891.625
$ws.Cells.Item(99, 11).Value = 8425.5
$ws.Cells.Item(99, 12).Value = 23674.875
$ws.Cells.Item(99, 13).Value = -6927.5
$ws.Cells.Item(99, 14).Value = -26670.875

# Row 137: Cutting Edge of Culinary Quality (item id 44013)
$ws.Cells.Item(137, 8).Value = 3926.197
$ws.Cells.Item(137, 9).Value = 4444.089
$ws.Cells.Item(137, 10).Value = 2816.4285
$ws.Cells.Item(137, 11).Value = 13332.267
$ws.Cells.Item(137, 12).Value = 8449.2855
$ws.Cells.Item(137, 13).Value = -10782.267
$ws.Cells.Item(137, 14).Value = -13549.2855

# Row 138: All-night Crafting (item id 44169)
$ws.Cells.Item(138, 8).Value = 5181.6787
$ws.Cells.Item(138, 9).Value = 1975.5
$ws.Cells.Item(138, 10).Value = 5428.3076
$ws.Cells.Item(138, 11).Value = 5926.5
$ws.Cells.Item(138, 12).Value = 16284.9228
$ws.Cells.Item(138, 13).Value = -786.5
$ws.Cells.Item(138, 14).Value = -26564.9228

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots (item id 27713)
$ws.Cells.Item(2, 8).Value = 4863.0557
$ws.Cells.Item(2, 9).Value = 4402.5
$ws.Cells.Item(2, 10).Value = 6475
$ws.Cells.Item(2, 11).Value = 4402.5
$ws.Cells.Item(2, 12).Value = 6475
$ws.Cells.Item(2, 13).Value = -4289.5
$ws.Cells.Item(2, 14).Value = -6701

# Row 32: Ingot We Trust (item id 44147)
$ws.Cells.Item(32, 8).Value = 2661.9365
$ws.Cells.Item(32, 9).Value = 1889.7192
$ws.Cells.Item(32, 10).Value = 9998
$ws.Cells.Item(32, 11).Value = 1889.7192
$ws.Cells.Item(32, 12).Value = 9998
$ws.Cells.Item(32, 13).Value = -1602.7192
$ws.Cells.Item(32, 14).Value = -10572

# Row 45: Hollow Hallmarks (item id 27714)
$ws.Cells.Item(45, 8).Value = 15490.607
$ws.Cells.Item(45, 9).Value = 12946.792
$ws.Cells.Item(45, 10).Value = 30753.5
$ws.Cells.Item(45, 11).Value = 12946.792
$ws.Cells.Item(45, 12).Value = 30753.5
$ws.Cells.Item(45, 13).Value = -12569.792
$ws.Cells.Item(45, 14).Value = -31507.5

# Row 61: Dealing with the Tough Stuff (item id 43999)
$ws.Cells.Item(61, 8).Value = 3136.3333
$ws.Cells.Item(61, 9).Value = 2912.087
$ws.Cells.Item(61, 10).Value = 4425.75
$ws.Cells.Item(61, 11).Value = 2912.087
$ws.Cells.Item(61, 12).Value = 4425.75
$ws.Cells.Item(61, 13).Value = -2700.087
$ws.Cells.Item(61, 14).Value = -4849.75

# Row 116: No Scope (item id 27713)
$ws.Cells.Item(116, 8).Value = 4863.0557
$ws.Cells.Item(116, 9).Value = 4402.5
$ws.Cells.Item(116, 10).Value = 6475
$ws.Cells.Item(116, 11).Value = 4402.5
$ws.Cells.Item(116, 12).Value = 6475
$ws.Cells.Item(116, 13).Value = -2108.5
$ws.Cells.Item(116, 14).Value = -11063

# Row 132: Don't Bore Me, Ore Me (item id 43997)
$ws.Cells.Item(132, 8).Value = 3612.1562
$ws.Cells.Item(132, 9).Value = 3061.1924
$ws.Cells.Item(132, 10).Value = 5999.6665
$ws.Cells.Item(132, 11).Value = 9183.5772
$ws.Cells.Item(132, 12).Value = 17998.9995
$ws.Cells.Item(132, 13).Value = -6653.5772
$ws.Cells.Item(132, 14).Value = -23058.9995

# Row 136: Metal with Mettle (item id 43999)
$ws.Cells.Item(136, 8).Value = 3136.3333
$ws.Cells.Item(136, 9).Value = 2912.087
$ws.Cells.Item(136, 10).Value = 4425.75
$ws.Cells.Item(136, 11).Value = 8736.261
$ws.Cells.Item(136, 12).Value = 13277.25
$ws.Cells.Item(136, 13).Value = -6186.261
$ws.Cells.Item(136, 14).Value = -18377.25

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells (item id 27713)
$ws.Cells.Item(3, 8).Value = 4863.0557
$ws.Cells.Item(3, 9).Value = 4402.5
$ws.Cells.Item(3, 10).Value = 6475
$ws.Cells.Item(3, 11).Value = 4402.5
$ws.Cells.Item(3, 12).Value = 6475
$ws.Cells.Item(3, 13).Value = -4288.5
$ws.Cells.Item(3, 14).Value = -6703

# Row 20: Smelt and Dealt (item id 14149)
$ws.Cells.Item(20, 8).Value = 3167.6428
$ws.Cells.Item(20, 9).Value = 2242.9697
$ws.Cells.Item(20, 10).Value = 6558.1113
$ws.Cells.Item(20, 11).Value = 2242.9697
$ws.Cells.Item(20, 12).Value = 6558.1113
$ws.Cells.Item(20, 13).Value = -1995.9697
$ws.Cells.Item(20, 14).Value = -7052.1113

# Row 22: Riveting Run (item id 5092)
$ws.Cells.Item(22, 8).Value = 281.33334
$ws.Cells.Item(22, 9).Value = 281.33334
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 281.33334
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -108.33334

# Row 36: I Saw What You Did There (item id 2320)
$ws.Cells.Item(36, 8).Value = 605.75
$ws.Cells.Item(36, 9).Value = 605.75
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 605.75
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -71.75

# Row 86: Through Thick and Thin (item id 12526)
$ws.Cells.Item(86, 8).Value = 5790.909
$ws.Cells.Item(86, 9).Value = 6133.1333
$ws.Cells.Item(86, 10).Value = 5057.5713
$ws.Cells.Item(86, 11).Value = 6133.1333
$ws.Cells.Item(86, 12).Value = 5057.5713
$ws.Cells.Item(86, 13).Value = -5010.1333
$ws.Cells.Item(86, 14).Value = -7303.5713

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) (item id 12526)
$ws.Cells.Item(89, 8).Value = 5790.909
$ws.Cells.Item(89, 9).Value = 6133.1333
$ws.Cells.Item(89, 10).Value = 5057.5713
$ws.Cells.Item(89, 11).Value = 30665.6665
$ws.Cells.Item(89, 12).Value = 25287.8565
$ws.Cells.Item(89, 13).Value = -25049.6665
$ws.Cells.Item(89, 14).Value = -36519.85649999999

# Row 107: The Gold Experience (item id 27706)
$ws.Cells.Item(107, 8).Value = 3570.25
$ws.Cells.Item(107, 9).Value = 3140.5
$ws.Cells.Item(107, 10).Value = 4000
$ws.Cells.Item(107, 11).Value = 3140.5
$ws.Cells.Item(107, 12).Value = 4000
$ws.Cells.Item(107, 13).Value = -1220.5
$ws.Cells.Item(107, 14).Value = -7840

# Row 134: Ruthenium Supremium (item id 43998)
$ws.Cells.Item(134, 8).Value = 3292.8667
$ws.Cells.Item(134, 9).Value = 2866.6667
$ws.Cells.Item(134, 10).Value = 4997.6665
$ws.Cells.Item(134, 11).Value = 8600.000100000001
$ws.Cells.Item(134, 12).Value = 14992.9995
$ws.Cells.Item(134, 13).Value = -6065.000100000001
$ws.Cells.Item(134, 14).Value = -20062.9995

# Row 139: Maul Me (item id 43261)
$ws.Cells.Item(139, 8).Value = 61022
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 61022
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 61022
$ws.Cells.Item(139, 14).Value = -71302

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found (item id 44023)
$ws.Cells.Item(31, 8).Value = 1865454
$ws.Cells.Item(31, 9).Value = 2534050
$ws.Cells.Item(31, 10).Value = 8242.944
$ws.Cells.Item(31, 11).Value = 2534050
$ws.Cells.Item(31, 12).Value = 8242.944
$ws.Cells.Item(31, 13).Value = -2533755
$ws.Cells.Item(31, 14).Value = -8832.944

# Row 34: Armoires of the Rich and Famous (item id 44023)
$ws.Cells.Item(34, 8).Value = 1865454
$ws.Cells.Item(34, 9).Value = 2534050
$ws.Cells.Item(34, 10).Value = 8242.944
$ws.Cells.Item(34, 11).Value = 2534050
$ws.Cells.Item(34, 12).Value = 8242.944
$ws.Cells.Item(34, 13).Value = -2533848
$ws.Cells.Item(34, 14).Value = -8646.944

# Row 132: Hull Lotta Damage (item id 44019)
$ws.Cells.Item(132, 8).Value = 3117.8125
$ws.Cells.Item(132, 9).Value = 2992.5806
$ws.Cells.Item(132, 10).Value = 7000
$ws.Cells.Item(132, 11).Value = 8977.7418
$ws.Cells.Item(132, 12).Value = 21000
$ws.Cells.Item(132, 13).Value = -6447.7418
$ws.Cells.Item(132, 14).Value = -26060

# Row 134: Wood You Be Quiet (item id 44020)
$ws.Cells.Item(134, 8).Value = 23039.266
$ws.Cells.Item(134, 9).Value = 19798.031
$ws.Cells.Item(134, 10).Value = 130000
$ws.Cells.Item(134, 11).Value = 59394.09299999999
$ws.Cells.Item(134, 12).Value = 390000
$ws.Cells.Item(134, 13).Value = -56859.09299999999
$ws.Cells.Item(134, 14).Value = -395070

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food (item id 4847)
$ws.Cells.Item(2, 8).Value = 90.333336
$ws.Cells.Item(2, 9).Value = 83.411766
$ws.Cells.Item(2, 10).Value = 102.1
$ws.Cells.Item(2, 11).Value = 500.470596
$ws.Cells.Item(2, 12).Value = 612.5999999999999
$ws.Cells.Item(2, 13).Value = -387.470596
$ws.Cells.Item(2, 14).Value = -838.5999999999999

# Row 12: Butter Me Up (item id 4854)
$ws.Cells.Item(12, 8).Value = 89.57895000000001
$ws.Cells.Item(12, 9).Value = 131
$ws.Cells.Item(12, 10).Value = 65.416664
$ws.Cells.Item(12, 11).Value = 393
$ws.Cells.Item(12, 12).Value = 196.249992
$ws.Cells.Item(12, 13).Value = -220
$ws.Cells.Item(12, 14).Value = -542.249992

# Row 38: Pretty as a Picture (item id 4860)
$ws.Cells.Item(38, 8).Value = 85
$ws.Cells.Item(38, 9).Value = 90
$ws.Cells.Item(38, 10).Value = 71
$ws.Cells.Item(38, 11).Value = 270
$ws.Cells.Item(38, 12).Value = 213
$ws.Cells.Item(38, 13).Value = 77
$ws.Cells.Item(38, 14).Value = -907

# Row 56: Culture Club (item id 10146)
$ws.Cells.Item(56, 8).Value = 10264.833
$ws.Cells.Item(56, 9).Value = 10264.833
$ws.Cells.Item(56, 10).Value = 0
$ws.Cells.Item(56, 11).Value = 10264.833
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(56, 13).Value = -9734.833000000001

# Row 86: Let's Not Get Sappy (item id 12892)
$ws.Cells.Item(86, 8).Value = 584
$ws.Cells.Item(86, 9).Value = 87
$ws.Cells.Item(86, 10).Value = 749.6667
$ws.Cells.Item(86, 11).Value = 261
$ws.Cells.Item(86, 12).Value = 2249.0001
$ws.Cells.Item(86, 13).Value = 925
$ws.Cells.Item(86, 14).Value = -4621.0001

# Row 89: Luxury Spillover (L) (item id 12892)
$ws.Cells.Item(89, 8).Value = 584
$ws.Cells.Item(89, 9).Value = 87
$ws.Cells.Item(89, 10).Value = 749.6667
$ws.Cells.Item(89, 11).Value = 783
$ws.Cells.Item(89, 12).Value = 6747.0003
$ws.Cells.Item(89, 13).Value = 5145
$ws.Cells.Item(89, 14).Value = -18603.0003

# Row 92: Oh No Udon (item id 19841)
$ws.Cells.Item(92, 8).Value = 1096.875
$ws.Cells.Item(92, 9).Value = 750
$ws.Cells.Item(92, 10).Value = 1146.4286
$ws.Cells.Item(92, 11).Value = 2250
$ws.Cells.Item(92, 12).Value = 3439.2858
$ws.Cells.Item(92, 13).Value = -1002
$ws.Cells.Item(92, 14).Value = -5935.2858

# Row 107: Slippery Service (item id 27838)
$ws.Cells.Item(107, 8).Value = 22223140
$ws.Cells.Item(107, 9).Value = 47619524
$ws.Cells.Item(107, 10).Value = 1301.25
$ws.Cells.Item(107, 11).Value = 142858572
$ws.Cells.Item(107, 12).Value = 3903.75
$ws.Cells.Item(107, 13).Value = -142856652
$ws.Cells.Item(107, 14).Value = -7743.75

# Row 124: Bobbing for Compliments (item id 36040)
$ws.Cells.Item(124, 8).Value = 20458.676
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 20458.676
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 61376.028
$ws.Cells.Item(124, 14).Value = -71196.02799999999
$ws.Cells.Item(124, 13).ClearContents()

# Row 129: Comfort Food (item id 36054)
$ws.Cells.Item(129, 8).Value = 5825838
$ws.Cells.Item(129, 9).Value = 16500758
$ws.Cells.Item(129, 10).Value = 3154.3635
$ws.Cells.Item(129, 11).Value = 49502274
$ws.Cells.Item(129, 12).Value = 9463.0905
$ws.Cells.Item(129, 13).Value = -49497274
$ws.Cells.Item(129, 14).Value = -19463.0905

# Row 131: The Mountain Steeped (item id 36060)
$ws.Cells.Item(131, 8).Value = 20383.834
$ws.Cells.Item(131, 9).Value = 101391.4
$ws.Cells.Item(131, 10).Value = 1973.0227
$ws.Cells.Item(131, 11).Value = 304174.2
$ws.Cells.Item(131, 12).Value = 5919.0681
$ws.Cells.Item(131, 13).Value = -299134.2
$ws.Cells.Item(131, 14).Value = -15999.0681

# Row 137: Creative Chocolate (item id 44088)
$ws.Cells.Item(137, 8).Value = 5142.75
$ws.Cells.Item(137, 9).Value = 5412.6665
$ws.Cells.Item(137, 10).Value = 4333
$ws.Cells.Item(137, 11).Value = 16237.9995
$ws.Cells.Item(137, 12).Value = 12999
$ws.Cells.Item(137, 13).Value = -11137.9995
$ws.Cells.Item(137, 14).Value = -23199

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence (item id 36182)
$ws.Cells.Item(122, 8).Value = 1671.5333
$ws.Cells.Item(122, 9).Value = 1527.3
$ws.Cells.Item(122, 10).Value = 1960
$ws.Cells.Item(122, 11).Value = 4581.9
$ws.Cells.Item(122, 12).Value = 5880
$ws.Cells.Item(122, 13).Value = -2131.9
$ws.Cells.Item(122, 14).Value = -10780

# Row 132: On Board for Lar (item id 44008)
$ws.Cells.Item(132, 8).Value = 2441.0908
$ws.Cells.Item(132, 9).Value = 2435.3125
$ws.Cells.Item(132, 10).Value = 2480.7144
$ws.Cells.Item(132, 11).Value = 7305.9375
$ws.Cells.Item(132, 12).Value = 7442.1432
$ws.Cells.Item(132, 13).Value = -4775.9375
$ws.Cells.Item(132, 14).Value = -12502.1432

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs (item id 5277)
$ws.Cells.Item(22, 8).Value = 1966.6666
$ws.Cells.Item(22, 9).Value = 2200
$ws.Cells.Item(22, 10).Value = 1500
$ws.Cells.Item(22, 11).Value = 2200
$ws.Cells.Item(22, 12).Value = 1500
$ws.Cells.Item(22, 13).Value = -1905
$ws.Cells.Item(22, 14).Value = -2090

# Row 27: Fire and Hide (item id 5277)
$ws.Cells.Item(27, 8).Value = 1966.6666
$ws.Cells.Item(27, 9).Value = 2200
$ws.Cells.Item(27, 10).Value = 1500
$ws.Cells.Item(27, 11).Value = 2200
$ws.Cells.Item(27, 12).Value = 1500
$ws.Cells.Item(27, 13).Value = -2093
$ws.Cells.Item(27, 14).Value = -1714

# Row 55: It's Not a Job, It's a Calling (item id 5284)
$ws.Cells.Item(55, 8).Value = 256.875
$ws.Cells.Item(55, 9).Value = 236.16667
$ws.Cells.Item(55, 10).Value = 319
$ws.Cells.Item(55, 11).Value = 236.16667
$ws.Cells.Item(55, 12).Value = 319
$ws.Cells.Item(55, 13).Value = -63.16667000000001
$ws.Cells.Item(55, 14).Value = -665

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire (item id 36208)
$ws.Cells.Item(122, 8).Value = 6156.684
$ws.Cells.Item(122, 9).Value = 3006.6365
$ws.Cells.Item(122, 10).Value = 10488
$ws.Cells.Item(122, 11).Value = 9019.9095
$ws.Cells.Item(122, 12).Value = 31464
$ws.Cells.Item(122, 13).Value = -6569.9095
$ws.Cells.Item(122, 14).Value = -36364

# Row 132: Comfy Cabins (item id 44029)
$ws.Cells.Item(132, 8).Value = 1922.0769
$ws.Cells.Item(132, 9).Value = 1972
$ws.Cells.Item(132, 10).Value = 1647.5
$ws.Cells.Item(132, 11).Value = 5916
$ws.Cells.Item(132, 12).Value = 4942.5
$ws.Cells.Item(132, 13).Value = -3386
$ws.Cells.Item(132, 14).Value = -10002.5

# Row 136: Weaving the Envelope (item id 44031)
$ws.Cells.Item(136, 8).Value = 3354.1365
$ws.Cells.Item(136, 9).Value = 3143.3125
$ws.Cells.Item(136, 10).Value = 3916.3333
$ws.Cells.Item(136, 11).Value = 9429.9375
$ws.Cells.Item(136, 12).Value = 11748.9999
$ws.Cells.Item(136, 13).Value = -6879.9375
$ws.Cells.Item(136, 14).Value = -16848.9999
